$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold text-formatted numbers (e.g. "71.087.62")
# and percentages padded with spaces. Force each touched cell to Text
# format before writing so Excel does not reinterpret the string as a
# number (which would both change its value and drop formatting like
# "10.60" -> 10.6).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.896.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.832.67'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '701.68'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.48'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.830.26'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.33'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.73'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.482.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.867.38'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '71.012.20'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.22'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.41'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '496.34'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.67'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.734'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.20'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000144'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.60'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.15'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.09'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.45'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.23'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.33'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.178'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.19'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.799.59'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.993'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.102'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.97'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.34'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.74%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.000316'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '163.62'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '431.75'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.95'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.72'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.38'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.297'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.16%  '
